$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(9, 8).Value = 92.933334
$ws.Cells.Item(9, 9).Value = 89.454544
$ws.Cells.Item(9, 10).Value = 102.5
$ws.Cells.Item(9, 11).Value = 89.454544
$ws.Cells.Item(9, 12).Value = 102.5
$ws.Cells.Item(9, 13).Value = 79.545456
$ws.Cells.Item(9, 14).Value = -440.5

$ws.Cells.Item(21, 8).Value = 1500
$ws.Cells.Item(21, 9).Value = 1500
$ws.Cells.Item(21, 11).Value = 1500
$ws.Cells.Item(21, 13).Value = -1032

$ws.Cells.Item(23, 8).Value = 1500
$ws.Cells.Item(23, 9).Value = 1500
$ws.Cells.Item(23, 11).Value = 1500
$ws.Cells.Item(23, 13).Value = -1266

$ws.Cells.Item(40, 13).ClearContents()
$ws.Cells.Item(40, 8).Value = 1745.4546
$ws.Cells.Item(40, 9).Value = 0
$ws.Cells.Item(40, 10).Value = 1745.4546
$ws.Cells.Item(40, 11).Value = 0
$ws.Cells.Item(40, 12).Value = 1745.4546
$ws.Cells.Item(40, 14).Value = -2095.4546

$ws.Cells.Item(51, 8).Value = 3720.32
$ws.Cells.Item(51, 9).Value = 4947.7334
$ws.Cells.Item(51, 10).Value = 1879.2
$ws.Cells.Item(51, 11).Value = 4947.7334
$ws.Cells.Item(51, 12).Value = 1879.2
$ws.Cells.Item(51, 13).Value = -4463.7334
$ws.Cells.Item(51, 14).Value = -2847.2

$ws.Cells.Item(132, 8).Value = 3863656.2
$ws.Cells.Item(132, 9).Value = 4610705.5
$ws.Cells.Item(132, 10).Value = 3902.6667
$ws.Cells.Item(132, 11).Value = 13832116.5
$ws.Cells.Item(132, 12).Value = 11708.0001
$ws.Cells.Item(132, 13).Value = -13829586.5
$ws.Cells.Item(132, 14).Value = -16768.0001

$ws.Cells.Item(137, 8).Value = 1192.9333
$ws.Cells.Item(137, 9).Value = 898.5
$ws.Cells.Item(137, 10).Value = 1300
$ws.Cells.Item(137, 11).Value = 2695.5
$ws.Cells.Item(137, 12).Value = 3900
$ws.Cells.Item(137, 13).Value = -145.5
$ws.Cells.Item(137, 14).Value = -9000

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(74, 8).Value = 664.5
$ws.Cells.Item(74, 9).Value = 688.5454999999999
$ws.Cells.Item(74, 10).Value = 400
$ws.Cells.Item(74, 11).Value = 688.5454999999999
$ws.Cells.Item(74, 12).Value = 400
$ws.Cells.Item(74, 13).Value = 185.4545000000001
$ws.Cells.Item(74, 14).Value = -2148

$ws.Cells.Item(77, 8).Value = 664.5
$ws.Cells.Item(77, 9).Value = 688.5454999999999
$ws.Cells.Item(77, 10).Value = 400
$ws.Cells.Item(77, 11).Value = 3442.7275
$ws.Cells.Item(77, 12).Value = 2000
$ws.Cells.Item(77, 13).Value = 925.2725
$ws.Cells.Item(77, 14).Value = -10736

$ws.Cells.Item(122, 8).Value = 1805.9584
$ws.Cells.Item(122, 9).Value = 1684.9445
$ws.Cells.Item(122, 11).Value = 5054.833500000001
$ws.Cells.Item(122, 13).Value = -2604.833500000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(86, 8).Value = 1700.3334
$ws.Cells.Item(86, 9).Value = 1233.8334
$ws.Cells.Item(86, 11).Value = 1233.8334
$ws.Cells.Item(86, 13).Value = -110.8334

$ws.Cells.Item(89, 8).Value = 1700.3334
$ws.Cells.Item(89, 9).Value = 1233.8334
$ws.Cells.Item(89, 11).Value = 6169.166999999999
$ws.Cells.Item(89, 13).Value = -553.1669999999995

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(7, 8).Value = 89.75
$ws.Cells.Item(7, 9).Value = 67.7
$ws.Cells.Item(7, 10).Value = 200
$ws.Cells.Item(7, 11).Value = 67.7
$ws.Cells.Item(7, 12).Value = 200
$ws.Cells.Item(7, 13).Value = 45.3
$ws.Cells.Item(7, 14).Value = -426

$ws.Cells.Item(16, 8).Value = 2037.8572
$ws.Cells.Item(16, 9).Value = 1237
$ws.Cells.Item(16, 10).Value = 2638.5
$ws.Cells.Item(16, 11).Value = 1237
$ws.Cells.Item(16, 12).Value = 2638.5
$ws.Cells.Item(16, 13).Value = -950
$ws.Cells.Item(16, 14).Value = -3212.5

$ws.Cells.Item(113, 8).Value = 2037.8572
$ws.Cells.Item(113, 9).Value = 1237
$ws.Cells.Item(113, 10).Value = 2638.5
$ws.Cells.Item(113, 11).Value = 1237
$ws.Cells.Item(113, 12).Value = 2638.5
$ws.Cells.Item(113, 13).Value = 933
$ws.Cells.Item(113, 14).Value = -6978.5

$ws.Cells.Item(122, 8).Value = 1245.4
$ws.Cells.Item(122, 9).Value = 890.8
$ws.Cells.Item(122, 10).Value = 1600
$ws.Cells.Item(122, 11).Value = 2672.4
$ws.Cells.Item(122, 12).Value = 4800
$ws.Cells.Item(122, 13).Value = -222.3999999999996
$ws.Cells.Item(122, 14).Value = -9700

$ws.Cells.Item(132, 8).Value = 2229.0557
$ws.Cells.Item(132, 9).Value = 1577.9
$ws.Cells.Item(132, 10).Value = 3043
$ws.Cells.Item(132, 11).Value = 4733.700000000001
$ws.Cells.Item(132, 12).Value = 9129
$ws.Cells.Item(132, 13).Value = -2203.700000000001
$ws.Cells.Item(132, 14).Value = -14189

$ws.Cells.Item(134, 8).Value = 812.5
$ws.Cells.Item(134, 9).Value = 823.86365
$ws.Cells.Item(134, 11).Value = 2471.59095
$ws.Cells.Item(134, 13).Value = 63.40905000000021

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(2, 8).Value = 41.416668
$ws.Cells.Item(2, 9).Value = 54.875
$ws.Cells.Item(2, 10).Value = 14.5
$ws.Cells.Item(2, 11).Value = 329.25
$ws.Cells.Item(2, 12).Value = 87
$ws.Cells.Item(2, 13).Value = -216.25
$ws.Cells.Item(2, 14).Value = -313

$ws.Cells.Item(122, 8).Value = 785.8333
$ws.Cells.Item(122, 9).Value = 433.07693
$ws.Cells.Item(122, 10).Value = 1703
$ws.Cells.Item(122, 11).Value = 3897.69237
$ws.Cells.Item(122, 12).Value = 15327
$ws.Cells.Item(122, 13).Value = -1447.69237
$ws.Cells.Item(122, 14).Value = -20227

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(2, 8).Value = 46.47059
$ws.Cells.Item(2, 9).Value = 15.428572
$ws.Cells.Item(2, 10).Value = 68.2
$ws.Cells.Item(2, 11).Value = 15.428572
$ws.Cells.Item(2, 12).Value = 68.2
$ws.Cells.Item(2, 13).Value = 97.571428
$ws.Cells.Item(2, 14).Value = -294.2

$ws.Cells.Item(118, 8).Value = 19800
$ws.Cells.Item(118, 10).Value = 19800
$ws.Cells.Item(118, 12).Value = 19800
$ws.Cells.Item(118, 14).Value = -23114

$ws.Cells.Item(132, 8).Value = 26543.146
$ws.Cells.Item(132, 9).Value = 32189.848
$ws.Cells.Item(132, 10).Value = 3250.5
$ws.Cells.Item(132, 11).Value = 96569.54400000001
$ws.Cells.Item(132, 12).Value = 9751.5
$ws.Cells.Item(132, 13).Value = -94039.54400000001
$ws.Cells.Item(132, 14).Value = -14811.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(132, 8).Value = 10315.714
$ws.Cells.Item(132, 9).Value = 16043.733
$ws.Cells.Item(132, 10).Value = 3706.4614
$ws.Cells.Item(132, 11).Value = 48131.199
$ws.Cells.Item(132, 12).Value = 11119.3842
$ws.Cells.Item(132, 13).Value = -45601.199
$ws.Cells.Item(132, 14).Value = -16179.3842

$ws.Cells.Item(136, 8).Value = 4253.0527
$ws.Cells.Item(136, 9).Value = 5714.8184
$ws.Cells.Item(136, 10).Value = 2243.125
$ws.Cells.Item(136, 11).Value = 17144.4552
$ws.Cells.Item(136, 12).Value = 6729.375
$ws.Cells.Item(136, 13).Value = -14594.4552
$ws.Cells.Item(136, 14).Value = -11829.375

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(74, 8).Value = 12285.2
$ws.Cells.Item(74, 10).Value = 12285.2
$ws.Cells.Item(74, 12).Value = 12285.2
$ws.Cells.Item(74, 14).Value = -14157.2

$ws.Cells.Item(77, 8).Value = 12285.2
$ws.Cells.Item(77, 10).Value = 12285.2
$ws.Cells.Item(77, 12).Value = 36855.60000000001
$ws.Cells.Item(77, 14).Value = -46215.60000000001

$ws.Cells.Item(96, 13).ClearContents()
$ws.Cells.Item(96, 8).Value = 2500
$ws.Cells.Item(96, 9).Value = 0
$ws.Cells.Item(96, 10).Value = 2500
$ws.Cells.Item(96, 11).Value = 0
$ws.Cells.Item(96, 12).Value = 2500
$ws.Cells.Item(96, 14).Value = -5246

$ws.Cells.Item(136, 8).Value = 1042.7273
$ws.Cells.Item(136, 9).Value = 1262.5807
$ws.Cells.Item(136, 10).Value = 518.46155
$ws.Cells.Item(136, 11).Value = 3787.7421
$ws.Cells.Item(136, 12).Value = 1555.38465
$ws.Cells.Item(136, 13).Value = -1237.7421
$ws.Cells.Item(136, 14).Value = -6655.38465
